# Replace the static "NBLxxxx" product-code values in B3:B8 with a
# generated code formula ("PK_KN_" & a random 5-digit number), matching
# the "Fix thêm excel phụ kiện next" commit.
#
# B3 gets its own (non-shared) formula; B4:B8 share one formula group,
# exactly as produced by Excel when you fill a formula down from B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$formula = '="PK_KN_"&TEXT(RANDBETWEEN(0,99999),"00000")'

$ws.Range("B3").Formula = $formula
$ws.Range("B4:B8").Formula = $formula

# The old values were long strings with embedded line breaks that forced
# the wrapped row height (ht="27.6"); the new short codes fit on one
# line, so let the rows shrink back to the sheet's default height.
$ws.Rows("3:8").EntireRow.AutoFit()

# Move the active selection, as recorded in the saved sheet view.
[void]$ws.Range("B13").Select()

Write-Output "applied PK_KN_ formula to B3:B8"
